try {
$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh3 = $m.Shapes.Item(3)
$tf2 = $sh3.TextFrame2
$tr2 = $tf2.TextRange
$runs = $tr2.Runs()
$i = 0
foreach ($r in $runs) {
    $i++
    if ($i -eq 5) {
        Write-Output "Run5 before: [$($r.Text)]"
        $r.Text = "Fifth Level"
        Write-Output "Run5 after: [$($r.Text)]"
    }
}
} catch {
Write-Output "ERROR: $_"
}
